# Generate Report for Handoff
# Updates Priority and Latest Handoff Datetime for the rows that were
# (re)handed off, on both localized-sheet tables, and propagates the new
# handoff datetime to the Overview sheet's "Latest HO Xliff Generate Date"
# column for the affected files.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsOverview = $wb.Worksheets.Item("Overview")

# zh-cn table: rows 4-7 -> Priority (E) "low" -> "ht"; Latest Handoff Datetime (H) updated
$wsZhCn.Range("E4:E7").Value = "ht"
$wsZhCn.Range("H4").Value = "2016-09-05 08:42:32"
$wsZhCn.Range("H5").Value = "2016-09-05 08:42:32"
$wsZhCn.Range("H6").Value = "2016-09-05 08:42:32"
$wsZhCn.Range("H7").Value = "2016-09-05 08:42:32"

# de-de table: rows 4-7 -> Priority (E) "low" -> "ht"; Latest Handoff Datetime (H) updated
$wsDeDe.Range("E4:E7").Value = "ht"
$wsDeDe.Range("H4").Value = "2016-09-05 08:42:37"
$wsDeDe.Range("H5").Value = "2016-09-05 08:42:37"
$wsDeDe.Range("H6").Value = "2016-09-05 08:42:37"
$wsDeDe.Range("H7").Value = "2016-09-05 08:42:37"

# Overview table: rows 4-7 -> Latest HO Xliff Generate Date (G) mirrors de-de handoff datetime
$wsOverview.Range("G4").Value = "2016-09-05 08:42:37"
$wsOverview.Range("G5").Value = "2016-09-05 08:42:37"
$wsOverview.Range("G6").Value = "2016-09-05 08:42:37"
$wsOverview.Range("G7").Value = "2016-09-05 08:42:37"
